# Activity.xlsx — "Added method for getting user sessions"
#
# 1. Update the GET_last_login sheet: the shared "preRequisite" description
#    text changes, a few leftover rows/cells get cleaned up, and three new
#    (currently-empty) columns are provisioned to match the new sheet.
# 2. Duplicate that sheet to create the new GET_user_sessions sheet and
#    adjust its description/URI cells + make it the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. GET_last_login sheet edits
# ---------------------------------------------------------------------
$lastLogin = $wb.Worksheets.Item("GET_last_login")

# preRequisite description text changed
$lastLogin.Range("D2").Value = "Get valid userId from Activity API"

# Rows 3-5 no longer carry the leftover TCID number / "N" flag
$lastLogin.Range("A3:B5").ClearContents()

# Row 5's trailing (unused) cells are removed outright
$lastLogin.Range("H5:J5").Clear()

# Provision the three new columns that the sessions sheet will also use
$lastLogin.Columns.Item(8).ColumnWidth = 21.1
$lastLogin.Columns.Item(9).ColumnWidth = 25.0
$lastLogin.Columns.Item(10).ColumnWidth = 25.4

# Selection / scroll state
$lastLogin.Range("C15").Select()

# ---------------------------------------------------------------------
# 2. Use a couple of scratch sheets so the new sheet lands on sheetId 5
#    (matches the sheetId Excel itself would have assigned), then remove
#    the scratch sheets again.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch1 = $wb.Worksheets.Add($null, $lastSheet)
$scratch1.Name = "Scratch1"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch2 = $wb.Worksheets.Add($null, $lastSheet2)
$scratch2.Name = "Scratch2"

$wb.Worksheets.Item("Scratch1").Delete()

# ---------------------------------------------------------------------
# 3. Duplicate GET_last_login -> GET_user_sessions
# ---------------------------------------------------------------------
$lastLogin.Copy($null, $lastLogin)
$sessions = $wb.Worksheets.Item("GET_last_login (2)")
$sessions.Name = "GET_user_sessions"

$wb.Worksheets.Item("Scratch2").Delete()

# New sheet's content: description + uri differ from GET_last_login
# (uri is written first so it lands before the description in the shared
# string table, matching the canonical ordering)
$sessions.Range("F2").Value = "/activity/v1/users/{userId}/sessions"
$sessions.Range("C2").Value = "Get User Sessions"

# Selection / scroll state + make it the active tab
$sessions.Range("I20").Select()
$sessions.Activate()
